$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview sheet, columns E/F on rows 2-3 share this text)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) Column width adjustments (report columns widened to fit new content)
# ---------------------------------------------------------------------------
$wide = 29.9777050018311 - 0.8333333333333334   # -> stored col width ~= 29.98
$full = 40.0 - 0.8333333333333334                # -> stored col width == 40

$ws1.Columns.Item(5).ColumnWidth = $wide   # Overview!E (zh-cn status)
$ws1.Columns.Item(6).ColumnWidth = $wide   # Overview!F (de-de status)

$ws2.Columns.Item(3).ColumnWidth = $wide   # zh-cn!C (Status)
$ws2.Columns.Item(9).ColumnWidth = $full   # zh-cn!I (Latest Target File)
$ws2.Columns.Item(10).ColumnWidth = $full  # zh-cn!J (Latest Handback File)

$ws3.Columns.Item(3).ColumnWidth = $wide   # de-de!C (Status)
$ws3.Columns.Item(9).ColumnWidth = $full   # de-de!I (Latest Target File)
$ws3.Columns.Item(10).ColumnWidth = $full  # de-de!J (Latest Handback File)

# ---------------------------------------------------------------------------
# 3) Populate handback report data for zh-cn
# ---------------------------------------------------------------------------
$zhTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad3815b877ad278e704ed11f8469c35b34ccc21d/e2e/5696414f-a192-4f49-8dcb-a4d802c52282.md"
$zhTargetName = "5696414f-a192-4f49-8dcb-a4d802c52282.md"
$zhHandbackFile = "5696414f-a192-4f49-8dcb-a4d802c52282.4c1b2034056c394f90553929ae4b0176f28ccb43.zh-cn.xlf"
$zhHandbackDate = "2016-10-13 13:37:26"

$ws2.Hyperlinks.Add($ws2.Range("I2"), $zhTargetUrl, [Type]::Missing, [Type]::Missing, $zhTargetName)
$ws2.Range("J2").Value = $zhHandbackFile
$ws2.Range("K2").Value = $zhHandbackDate

$ws2.Hyperlinks.Add($ws2.Range("I3"), $zhTargetUrl, [Type]::Missing, [Type]::Missing, $zhTargetName)
$ws2.Range("J3").Value = $zhHandbackFile
$ws2.Range("K3").Value = $zhHandbackDate

# ---------------------------------------------------------------------------
# 4) Populate handback report data for de-de
# ---------------------------------------------------------------------------
$deTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad3815b877ad278e704ed11f8469c35b34ccc21d/e2e/5696414f-a192-4f49-8dcb-a4d802c52282.md"
$deTargetName = "5696414f-a192-4f49-8dcb-a4d802c52282.md"
$deHandbackFile = "5696414f-a192-4f49-8dcb-a4d802c52282.4c1b2034056c394f90553929ae4b0176f28ccb43.de-de.xlf"
$deHandbackDate = "2016-10-13 13:37:43"

$ws3.Hyperlinks.Add($ws3.Range("I2"), $deTargetUrl, [Type]::Missing, [Type]::Missing, $deTargetName)
$ws3.Range("J2").Value = $deHandbackFile
$ws3.Range("K2").Value = $deHandbackDate

$ws3.Hyperlinks.Add($ws3.Range("I3"), $deTargetUrl, [Type]::Missing, [Type]::Missing, $deTargetName)
$ws3.Range("J3").Value = $deHandbackFile
$ws3.Range("K3").Value = $deHandbackDate
